$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" '70.509.16'
Set-TextValue $ws "E2" '  +0.37%  '
Set-TextValue $ws "D3" '3.569.04'
Set-TextValue $ws "E3" '  -0.35%  '
Set-TextValue $ws "E4" '  -0.06%  '
Set-TextValue $ws "D5" '583.30'
Set-TextValue $ws "E5" '  +0.87%  '
Set-TextValue $ws "D6" '183.74'
Set-TextValue $ws "E6" '  -1.74%  '
Set-TextValue $ws "D7" '3.558.75'
Set-TextValue $ws "E7" '  -0.52%  '
Set-TextValue $ws "E8" '  -0.29%  '
Set-TextValue $ws "E9" '  +0.06%  '
Set-TextValue $ws "D10" '0.220'
Set-TextValue $ws "E10" '  +16.49%  '
Set-TextValue $ws "D11" '0.647'
Set-TextValue $ws "E11" '  -1.35%  '
Set-TextValue $ws "D12" '54.01'
Set-TextValue $ws "E12" '  -1.08%  '
Set-TextValue $ws "D13" '0.0000321'
Set-TextValue $ws "E13" '  +5.02%  '
Set-TextValue $ws "D14" '9.49'
Set-TextValue $ws "E14" '  -1.10%  '
Set-TextValue $ws "D15" '4.134.18'
Set-TextValue $ws "E15" '  -0.58%  '
Set-TextValue $ws "D16" '19.49'
Set-TextValue $ws "E16" '  -1.19%  '
Set-TextValue $ws "D17" '70.433.18'
Set-TextValue $ws "E17" '  +0.36%  '
Set-TextValue $ws "D18" '3.547.02'
Set-TextValue $ws "E18" '  -0.62%  '
Set-TextValue $ws "D19" '568.87'
Set-TextValue $ws "E19" '  +14.78%  '
Set-TextValue $ws "D20" '12.33'
Set-TextValue $ws "E20" '  -1.52%  '
Set-TextValue $ws "E21" '  -0.26%  '
Set-TextValue $ws "E22" '  -3.57%  '
Set-TextValue $ws "D23" '17.61'
Set-TextValue $ws "E23" '  -9.04%  '
Set-TextValue $ws "D24" '4.57'
Set-TextValue $ws "E24" '  +4.01%  '
Set-TextValue $ws "D25" '4.97'
Set-TextValue $ws "E25" '  -0.87%  '
Set-TextValue $ws "D26" '95.39'
Set-TextValue $ws "E26" '  -1.52%  '
Set-TextValue $ws "D27" '11.24'
Set-TextValue $ws "E27" '  -3.43%  '
Set-TextValue $ws "D28" '2.92'
Set-TextValue $ws "E28" '  -2.59%  '
Set-TextValue $ws "D29" '9.09'
Set-TextValue $ws "E29" '  -3.15%  '
Set-TextValue $ws "D30" '32.17'
Set-TextValue $ws "E30" '  +1.38%  '
Set-TextValue $ws "D31" '7.30'
Set-TextValue $ws "E31" '  -6.30%  '
Set-TextValue $ws "D32" '12.22'
Set-TextValue $ws "E32" '  -4.58%  '
Set-TextValue $ws "B33" 'OKB'
Set-TextValue $ws "C33" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws "D33" '64.17'
Set-TextValue $ws "E33" '  -2.22%  '
Set-TextValue $ws "B34" 'Hedera'
Set-TextValue $ws "C34" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D34" '0.114'
Set-TextValue $ws "E34" '  -1.58%  '
Set-TextValue $ws "D35" '3.31'
Set-TextValue $ws "E35" '  +0.35%  '
Set-TextValue $ws "D36" '553.04'
Set-TextValue $ws "E36" '  -3.82%  '
Set-TextValue $ws "D37" '0.411'
Set-TextValue $ws "E37" '  -0.32%  '
Set-TextValue $ws "D38" '0.0₃0812'
Set-TextValue $ws "E38" '  +2.14%  '
Set-TextValue $ws "E39" '  +0.00%  '
Set-TextValue $ws "D40" '37.47'
Set-TextValue $ws "E40" '  -4.03%  '
Set-TextValue $ws "B41" 'Maker'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws "D41" '3.407.19'
Set-TextValue $ws "E41" '  +5.14%  '
Set-TextValue $ws "B42" 'Kaspa'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws "D42" '0.135'
Set-TextValue $ws "E42" '  -0.08%  '
Set-TextValue $ws "B43" 'Stacks'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws "D43" '3.35'
Set-TextValue $ws "E43" '  -3.43%  '
Set-TextValue $ws "B44" 'dogwifhat'
Set-TextValue $ws "C44" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws "D44" '3.06'
Set-TextValue $ws "E44" '  -4.71%  '
Set-TextValue $ws "B45" 'ApeXProtocol'
Set-TextValue $ws "C45" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws "D45" '3.56'
Set-TextValue $ws "E45" '  -6.03%  '
Set-TextValue $ws "B46" 'VeChain'
Set-TextValue $ws "C46" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws "D46" '0.0442'
Set-TextValue $ws "E46" '  -2.78%  '
Set-TextValue $ws "D47" '2.95'
Set-TextValue $ws "E47" '  -3.55%  '
Set-TextValue $ws "D48" '9.42'
Set-TextValue $ws "E48" '  -0.34%  '
Set-TextValue $ws "E49" '  +1.19%  '
Set-TextValue $ws "E50" '  +0.18%  '
Set-TextValue $ws "D51" '1.44'
Set-TextValue $ws "E51" '  -3.99%  '
